$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-5 with new data results (columns B:E updated, F is the total)
$ws.Range("B2:E2").Value = 9
$ws.Range("F2").Value = 36

$ws.Range("B3:E3").Value = 9
$ws.Range("F3").Value = 36

$ws.Range("B4:E4").Value = 10
$ws.Range("F4").Value = 40

$ws.Range("B5:E5").Value = 7
$ws.Range("F5").Value = 28

# Remove row 10 entirely (data no longer present)
$ws.Rows.Item(10).Delete()
